$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cliente" column (E) used placeholder first-names; replace them with the
# full customer names used in the updated learning-guide resource list.
$clientes = @{
    "E6"  = "Cherokee Macias"
    "E7"  = "Vladimir Potter"
    "E8"  = "David Mckenzie"
    "E9"  = "Cassady Leblanc"
    "E12" = "Uriah Knapp"
    "E13" = "Emerald Black"
    "E14" = "David Mckenzie"
    "E15" = "Xantha Beck"
    "E17" = "Virginia Holman"
    "E18" = "Zenia Cameron"
    "E19" = "Cassady Leblanc"
}

foreach ($addr in $clientes.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $clientes[$addr]
    # Give the refreshed customer names their own (black) font so the
    # updated values are visually distinguished from the rest of the sheet.
    $cell.Font.ColorIndex = 1
}

# Restore the saved cursor position from the edited workbook.
$ws.Range("K14").Select()
